$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.665.26'
$ws.Range("E2").Value = '  +3.58%  '
$ws.Range("D3").Value = '1.603.35'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("E6").Value = '  -0.53%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '27.09'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +9.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.41'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.86%  '
$ws.Range("E10").Value = '  +1.82%  '
$ws.Range("E11").Value = '  +1.91%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0908'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.88%  '
$ws.Range("D13").Value = '1.832.18'
$ws.Range("E13").Value = '  +2.43%  '
$ws.Range("D14").Value = '1.605.50'
$ws.Range("E14").Value = '  +2.78%  '
$ws.Range("D15").Value = '29.653.92'
$ws.Range("E15").Value = '  +3.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.536'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.75%  '
$ws.Range("E17").Value = '  +2.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.35'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '241.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.14%  '
$ws.Range("E20").Value = '  +3.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.94%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.70'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.74%  '
$ws.Range("E27").Value = '  +3.50%  '
$ws.Range("E28").Value = '  +0.42%  '
$ws.Range("E29").Value = '  +2.39%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("E31").Value = '  +3.65%  '
$ws.Range("E32").Value = '  +0.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.21'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.71%  '
$ws.Range("D35").Value = '1.428.87'
$ws.Range("E35").Value = '  +1.74%  '
$ws.Range("E36").Value = '  +0.37%  '
$ws.Range("E37").Value = '  +4.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.79'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.82%  '
$ws.Range("E39").Value = '  +0.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0166'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.88%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.538'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '54.84'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +30.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.96'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.07%  '
$ws.Range("E44").Value = '  +6.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '65.90'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.944'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +13.35%  '
$ws.Range("E49").Value = '  +0.94%  '
$ws.Range("D50").Value = '1.742.75'
$ws.Range("E50").Value = '  +2.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '86.45'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.05%  '
